$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("On-Site")

# Align A3's formatting with the rest of the row (was using the blank
# sub-system column style, now matches the numeric cell style)
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate row 3 with the "Mooring Lines" sub-system spare parts data
$ws.Range("A3").Value = "Mooring Lines"
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3

# Make "On-Site" the active sheet with A3 selected
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
